# Auto-generated edit script: updates market-price derived columns (H-N)
# across ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 10782.267
$ws.Range("I62").Value = 13210.2
$ws.Range("J62").Value = 5926.4
$ws.Range("K62").Value = 13210.2
$ws.Range("L62").Value = 5926.4
$ws.Range("M62").Value = -12586.2
$ws.Range("N62").Value = -7174.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 10782.267
$ws.Range("I65").Value = 13210.2
$ws.Range("J65").Value = 5926.4
$ws.Range("K65").Value = 66051
$ws.Range("L65").Value = 29632
$ws.Range("M65").Value = -62931
$ws.Range("N65").Value = -35872

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3990.5925
$ws.Range("I74").Value = 3714.6155
$ws.Range("J74").Value = 4246.857
$ws.Range("K74").Value = 3714.6155
$ws.Range("L74").Value = 4246.857
$ws.Range("M74").Value = -2778.6155
$ws.Range("N74").Value = -6118.857

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3990.5925
$ws.Range("I77").Value = 3714.6155
$ws.Range("J77").Value = 4246.857
$ws.Range("K77").Value = 18573.0775
$ws.Range("L77").Value = 21234.285
$ws.Range("M77").Value = -13893.0775
$ws.Range("N77").Value = -30594.285

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3435.2856
$ws.Range("I116").Value = 3550
$ws.Range("J116").Value = 3282.3333
$ws.Range("K116").Value = 3550
$ws.Range("L116").Value = 3282.3333
$ws.Range("M116").Value = -108
$ws.Range("N116").Value = -10166.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1983.1632
$ws.Range("I132").Value = 1061.3334
$ws.Range("J132").Value = 7514.143
$ws.Range("K132").Value = 3184.0002
$ws.Range("L132").Value = 22542.429
$ws.Range("M132").Value = -654.0001999999999
$ws.Range("N132").Value = -27602.429

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 789
$ws.Range("I135").Value = 443
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 3987
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -1452
$ws.Range("N135").Value = -23070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 700
$ws.Range("I15").Value = 700
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 700
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -350
$ws.Range("N15").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 9000
$ws.Range("J18").Value = 9000
$ws.Range("L18").Value = 9000
$ws.Range("N18").Value = -9644

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 213235.81
$ws.Range("I74").Value = 278806.88
$ws.Range("K74").Value = 278806.88
$ws.Range("M74").Value = -277932.88

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 213235.81
$ws.Range("I77").Value = 278806.88
$ws.Range("K77").Value = 1394034.4
$ws.Range("M77").Value = -1389666.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1860.3096
$ws.Range("I134").Value = 1033.3214
$ws.Range("J134").Value = 3514.2856
$ws.Range("K134").Value = 3099.9642
$ws.Range("L134").Value = 10542.8568
$ws.Range("M134").Value = -564.9642000000003
$ws.Range("N134").Value = -15612.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2959.8206
$ws.Range("I31").Value = 1981.1
$ws.Range("J31").Value = 6222.222
$ws.Range("K31").Value = 1981.1
$ws.Range("L31").Value = 6222.222
$ws.Range("M31").Value = -1686.1
$ws.Range("N31").Value = -6812.222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2959.8206
$ws.Range("I34").Value = 1981.1
$ws.Range("J34").Value = 6222.222
$ws.Range("K34").Value = 1981.1
$ws.Range("L34").Value = 6222.222
$ws.Range("M34").Value = -1779.1
$ws.Range("N34").Value = -6626.222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 822.069
$ws.Range("I105").Value = 809.2308
$ws.Range("K105").Value = 809.2308
$ws.Range("M105").Value = 937.7692

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1287.7778
$ws.Range("I107").Value = 193
$ws.Range("J107").Value = 1506.7333
$ws.Range("K107").Value = 579
$ws.Range("L107").Value = 4520.199900000001
$ws.Range("M107").Value = 1341
$ws.Range("N107").Value = -8360.1999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2854.2307
$ws.Range("I7").Value = 2750
$ws.Range("J7").Value = 3021
$ws.Range("K7").Value = 2750
$ws.Range("L7").Value = 3021
$ws.Range("M7").Value = -2638
$ws.Range("N7").Value = -3245

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 13005.556
$ws.Range("I40").Value = 14443.75
$ws.Range("J40").Value = 1500
$ws.Range("K40").Value = 14443.75
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = -14307.75
$ws.Range("N40").Value = -1772

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1038.4166
$ws.Range("I46").Value = 1177.2222
$ws.Range("J46").Value = 622
$ws.Range("K46").Value = 1177.2222
$ws.Range("L46").Value = 622
$ws.Range("M46").Value = -989.2221999999999
$ws.Range("N46").Value = -998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1421.3334
$ws.Range("I61").Value = 1421.3334
$ws.Range("K61").Value = 1421.3334
$ws.Range("M61").Value = -1219.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 817.6539
$ws.Range("I93").Value = 750.3913
$ws.Range("J93").Value = 1333.3334
$ws.Range("K93").Value = 750.3913
$ws.Range("L93").Value = 1333.3334
$ws.Range("M93").Value = 497.6087
$ws.Range("N93").Value = -3829.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1421.3334
$ws.Range("I113").Value = 1421.3334
$ws.Range("K113").Value = 1421.3334
$ws.Range("M113").Value = 748.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2854.2307
$ws.Range("I126").Value = 2750
$ws.Range("J126").Value = 3021
$ws.Range("K126").Value = 8250
$ws.Range("L126").Value = 9063
$ws.Range("M126").Value = -5780
$ws.Range("N126").Value = -14003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8889.666999999999
$ws.Range("I132").Value = 2389
$ws.Range("J132").Value = 18640.666
$ws.Range("K132").Value = 7167
$ws.Range("L132").Value = 55921.99800000001
$ws.Range("M132").Value = -4637
$ws.Range("N132").Value = -60981.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4841.8857
$ws.Range("I136").Value = 2774.2354
$ws.Range("J136").Value = 6794.6665
$ws.Range("K136").Value = 8322.706200000001
$ws.Range("L136").Value = 20383.9995
$ws.Range("M136").Value = -5772.706200000001
$ws.Range("N136").Value = -25483.9995
